$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H110").Value = 9999
$ws.Range("J110").Value = 9999
$ws.Range("L110").Value = 9999
$ws.Range("N110").Value = -18179
$ws.Range("H116").Value = 5999.5
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H132").Value = 1894.7858
$ws.Range("I132").Value = 1039.0741
$ws.Range("K132").Value = 3117.2223
$ws.Range("M132").Value = -587.2223000000004
$ws.Range("H138").Value = 2505.3125
$ws.Range("J138").Value = 3372
$ws.Range("L138").Value = 10116
$ws.Range("N138").Value = -20396

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12320.903
$ws.Range("I32").Value = 9547.25
$ws.Range("K32").Value = 9547.25
$ws.Range("M32").Value = -9260.25
$ws.Range("H88").Value = 1913.6364
$ws.Range("I88").Value = 2170
$ws.Range("J88").Value = 1700
$ws.Range("K88").Value = 2170
$ws.Range("L88").Value = 1700
$ws.Range("M88").Value = -1764
$ws.Range("N88").Value = -2512
$ws.Range("H91").Value = 1913.6364
$ws.Range("I91").Value = 2170
$ws.Range("J91").Value = 1700
$ws.Range("K91").Value = 2170
$ws.Range("L91").Value = 1700
$ws.Range("M91").Value = -766
$ws.Range("N91").Value = -4508

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4504.25
$ws.Range("I86").Value = 3273.5
$ws.Range("J86").Value = 6719.6
$ws.Range("K86").Value = 3273.5
$ws.Range("L86").Value = 6719.6
$ws.Range("M86").Value = -2150.5
$ws.Range("N86").Value = -8965.6
$ws.Range("H89").Value = 4504.25
$ws.Range("I89").Value = 3273.5
$ws.Range("J89").Value = 6719.6
$ws.Range("K89").Value = 16367.5
$ws.Range("L89").Value = 33598
$ws.Range("M89").Value = -10751.5
$ws.Range("N89").Value = -44830
$ws.Range("H105").Value = 2440.3635
$ws.Range("I105").Value = 1913.4286
$ws.Range("K105").Value = 1913.4286
$ws.Range("M105").Value = -166.4286
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H116").Value = 64535
$ws.Range("J116").Value = 64535
$ws.Range("L116").Value = 64535
$ws.Range("N116").Value = -73713

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 2875.1667
$ws.Range("I6").Value = 500
$ws.Range("J6").Value = 29002
$ws.Range("K6").Value = 500
$ws.Range("L6").Value = 29002
$ws.Range("M6").Value = -387
$ws.Range("N6").Value = -29228
$ws.Range("H7").Value = 131
$ws.Range("I7").Value = 106.85714
$ws.Range("K7").Value = 106.85714
$ws.Range("M7").Value = 6.142859999999999
$ws.Range("H25").Value = 5735
$ws.Range("I25").Value = 4974.5
$ws.Range("J25").Value = 6495.5
$ws.Range("K25").Value = 4974.5
$ws.Range("L25").Value = 6495.5
$ws.Range("M25").Value = -4800.5
$ws.Range("N25").Value = -6843.5
$ws.Range("H31").Value = 5799.423
$ws.Range("I31").Value = 4826.05
$ws.Range("J31").Value = 9044
$ws.Range("K31").Value = 4826.05
$ws.Range("L31").Value = 9044
$ws.Range("M31").Value = -4531.05
$ws.Range("N31").Value = -9634
$ws.Range("H34").Value = 5799.423
$ws.Range("I34").Value = 4826.05
$ws.Range("J34").Value = 9044
$ws.Range("K34").Value = 4826.05
$ws.Range("L34").Value = 9044
$ws.Range("M34").Value = -4624.05
$ws.Range("N34").Value = -9448
$ws.Range("H58").Value = 3056.0833
$ws.Range("I58").Value = 1913.2632
$ws.Range("K58").Value = 1913.2632
$ws.Range("M58").Value = -1710.2632
$ws.Range("H99").Value = 5075.2666
$ws.Range("I99").Value = 4723.5713
$ws.Range("K99").Value = 4723.5713
$ws.Range("M99").Value = -3225.5713
$ws.Range("H126").Value = 5075.2666
$ws.Range("I126").Value = 4723.5713
$ws.Range("K126").Value = 14170.7139
$ws.Range("M126").Value = -11700.7139
$ws.Range("H134").Value = 2010.6571
$ws.Range("I134").Value = 1841.0588
$ws.Range("K134").Value = 5523.1764
$ws.Range("M134").Value = -2988.1764
$ws.Range("H136").Value = 3056.0833
$ws.Range("I136").Value = 1913.2632
$ws.Range("K136").Value = 5739.7896
$ws.Range("M136").Value = -3189.7896

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 43.75
$ws.Range("I2").Value = 42.5
$ws.Range("J2").Value = 46.25
$ws.Range("K2").Value = 255
$ws.Range("L2").Value = 277.5
$ws.Range("M2").Value = -142
$ws.Range("N2").Value = -503.5
$ws.Range("H34").Value = 2604.8
$ws.Range("J34").Value = 3678.5715
$ws.Range("L34").Value = 11035.7145
$ws.Range("N34").Value = -11203.7145
$ws.Range("H39").Value = 7574.75
$ws.Range("J39").Value = 9999.666999999999
$ws.Range("L39").Value = 29999.001
$ws.Range("N39").Value = -30587.001
$ws.Range("H131").Value = 998.375
$ws.Range("I131").Value = 998.1429000000001
$ws.Range("K131").Value = 2994.4287
$ws.Range("M131").Value = 2045.5713

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 97888.336
$ws.Range("I132").Value = 126073.664
$ws.Range("K132").Value = 378220.992
$ws.Range("M132").Value = -375690.992

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H54").Value = 13971.25
$ws.Range("J54").Value = 13971.25
$ws.Range("L54").Value = 13971.25
$ws.Range("N54").Value = -15259.25
$ws.Range("H93").Value = 1687
$ws.Range("I93").Value = 1533
$ws.Range("J93").Value = 1733.2
$ws.Range("K93").Value = 1533
$ws.Range("L93").Value = 1733.2
$ws.Range("M93").Value = -285
$ws.Range("N93").Value = -4229.2
$ws.Range("H132").Value = 7967.0625
$ws.Range("I132").Value = 5943.3335
$ws.Range("K132").Value = 17830.0005
$ws.Range("M132").Value = -15300.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 43292
$ws.Range("J41").Value = 44938
$ws.Range("L41").Value = 44938
$ws.Range("N41").Value = -45718
$ws.Range("H63").Value = 32999.668
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 32999.668
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 32999.668
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -34247.668
$ws.Range("H66").Value = 32999.668
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 32999.668
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 98999.00399999999
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -105239.004
$ws.Range("H68").Value = 31999
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 31999
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 31999
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -33621
$ws.Range("H71").Value = 31999
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 31999
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 95997
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -104109
$ws.Range("H81").Value = 3642.7144
$ws.Range("I81").Value = 3916.5
$ws.Range("K81").Value = 7833
$ws.Range("M81").Value = -6772
$ws.Range("H84").Value = 3642.7144
$ws.Range("I84").Value = 3916.5
$ws.Range("K84").Value = 39165
$ws.Range("M84").Value = -33861
$ws.Range("H122").Value = 3168.3333
$ws.Range("I122").Value = 2391.8235
$ws.Range("J122").Value = 5054.143
$ws.Range("K122").Value = 7175.470499999999
$ws.Range("L122").Value = 15162.429
$ws.Range("M122").Value = -4725.470499999999
$ws.Range("N122").Value = -20062.429
